$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 should look like the other header cells (bold, bordered, centered/top)
# Copy the formatting from the existing header cell E1 onto F1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the cell values
$ws.Range("F1").Value = "Modelo"
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
